$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell value changes (rows 2-25) ---
# Row 3: E3 was empty -> -5.7
$ws.Cells.Item(3, 5).Value = -5.7

# Row 4: F4 was 17.97 -> now empty
$ws.Cells.Item(4, 6).ClearContents()

# Row 5: E5 was -5 -> now empty
$ws.Cells.Item(5, 5).ClearContents()

# Row 9: F9 was empty -> 17.26
$ws.Cells.Item(9, 6).Value = 17.26

# Row 10: F10 was empty -> 16.43
$ws.Cells.Item(10, 6).Value = 16.43

# Row 13: F13 was 17.1 -> now empty
$ws.Cells.Item(13, 6).ClearContents()

# Row 14: F14 was 17.76 -> now empty
$ws.Cells.Item(14, 6).ClearContents()

# Row 21: E21 was empty -> -8.699999999999999
$ws.Cells.Item(21, 5).Value = -8.699999999999999

# Row 23: E23 was -7 -> now empty
$ws.Cells.Item(23, 5).ClearContents()

# --- Remove two whole records (rows) ---
# "RM 232" (row 26) is removed entirely; all rows below shift up by one.
$ws.Rows.Item(26).Delete()

# After the shift above, "SC 92" (originally row 28) is now at row 27;
# remove it as well so every following row shifts up again.
$ws.Rows.Item(27).Delete()

# "SC 193" (originally row 34) is now at row 32 after both deletions above;
# its E value was empty and is now -6.4.
$ws.Cells.Item(32, 5).Value = -6.4
